{"js": "// \"Correccion ciudad en certificados\"\n// Replace the {{facultad.ciudad}} merge field with the literal, hard-coded\n// city name \"San Rafael, Mendoza\" and drop the now-redundant leading\n// period before \" el \" (the sentence used to read\n// \"... en {{facultad.ciudad}}. el {{fecha}}.-\" and now reads\n// \"... en San Rafael, Mendoza el {{fecha}}.-\").\n\n// Step 1: swap the {{facultad.ciudad}} field for the literal city text,\n// keeping the bold run formatting that wrapped the original field.\nconst cityResults = context.document.body.search(\"{{facultad.ciudad}}\", { matchCase: true, matchWildcards: false });\ncityResults.load(\"items\");\nawait context.sync();\n\nif (cityResults.items.length > 0) {\n  cityResults.items[0].insertText(\"San Rafael, Mendoza \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Step 2: the sentence used to have \". el \" right after the field; now\n// that the field text already ends with a space (\"Mendoza \"), the\n// leading period/space pair collapses to just \"el \".\nconst periodResults = context.document.body.search(\". el \", { matchCase: true, matchWildcards: false });\nperiodResults.load(\"items\");\nawait context.sync();\n\nif (periodResults.items.length > 0) {\n  periodResults.items[0].insertText(\"el \", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Correccion ciudad en certificados\"\n# Replace the {{facultad.ciudad}} merge field with the literal, hard-coded\n# city name \"San Rafael, Mendoza\" and drop the now-redundant leading\n# period before \" el \" (the sentence used to read\n# \"... en {{facultad.ciudad}}. el {{fecha}}.-\" and now reads\n# \"... en San Rafael, Mendoza el {{fecha}}.-\").\n\n$d = $word.ActiveDocument\n\n# Step 1: swap the {{facultad.ciudad}} field for the literal city text,\n# keeping the bold run formatting that wrapped the original field.\n$cityRange = $d.Content\n$cityFind = $cityRange.Find\n$cityFind.Text = \"{{facultad.ciudad}}\"\n$cityFind.MatchWildcards = $false\n$cityFound = $cityFind.Execute()\nif ($cityFound) {\n    $cityRange.Text = \"San Rafael, Mendoza \"\n}\n\n# Step 2: the sentence used to have \". el \" right after the field; now\n# that the field text already ends with a space (\"Mendoza \"), the\n# leading period/space pair collapses to just \"el \".\n$periodRange = $d.Content\n$periodFind = $periodRange.Find\n$periodFind.Text = \". el \"\n$periodFind.MatchWildcards = $false\n$periodFound = $periodFind.Execute()\nif ($periodFound) {\n    $periodRange.Text = \"el \"\n}\n"}
